# feat: add 2022-Q1 data
#
# The workbook currently has two sheets: "2021-Q4" and "总计".
# This inserts a new "2022-Q1" sheet (as a copy of "2021-Q4", so it
# inherits the same layout/styles) right after "2021-Q4", fills it with
# the 2022-Q1 holdings data, and updates the "总计" summary sheet with a
# new first data row for 2022-Q1 (pushing the existing 2021-Q4 row down).

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" sheet right after "2021-Q4" ---------------
$template.Copy($null, $template)
$q1 = $wb.Worksheets.Item(2)
$q1.Name = "2022-Q1"

# Look this up only now (by name) - inserting the new sheet shifts the
# "总计" sheet from index 2 to index 3, so a reference fetched beforehand
# would end up pointing at the wrong sheet.
$summary = $wb.Worksheets.Item("总计")

# --- 2. Fill in the 2022-Q1 holdings data -----------------------------------
# Columns D,E,F,G and the fund-code column B hold numeric-looking text
# (e.g. leading zeros / trailing zeros that must be preserved), so force
# text storage via NumberFormat "@" before assigning the value.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
Set-TextValue $q1.Range("B2") "007178"
$q1.Range("C2").Value = "浙商港股通中华交易服务预期高股息指数增强A"
Set-TextValue $q1.Range("D2") "7.93"
Set-TextValue $q1.Range("E2") "90.20"
Set-TextValue $q1.Range("F2") "7.73"
Set-TextValue $q1.Range("G2") "0.6130"
$q1.Range("H2").Value = 3

# Row 3
Set-TextValue $q1.Range("B3") "007216"
$q1.Range("C3").Value = "浙商港股通中华交易服务预期高股息指数增强C"
Set-TextValue $q1.Range("D3") "4.60"
Set-TextValue $q1.Range("E3") "90.20"
Set-TextValue $q1.Range("F3") "7.73"
Set-TextValue $q1.Range("G3") "0.3556"
$q1.Range("H3").Value = 3

# Row 4 (new row)
$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "011018"
$q1.Range("C4").Value = "景顺长城安泽回报一年持有期混合A"
Set-TextValue $q1.Range("D4") "7.81"
Set-TextValue $q1.Range("E4") "34.82"
Set-TextValue $q1.Range("F4") "0.70"
Set-TextValue $q1.Range("G4") "0.0547"
$q1.Range("H4").Value = 10

# Row 5 (new row)
$q1.Range("A5").Value = 3
Set-TextValue $q1.Range("B5") "011019"
$q1.Range("C5").Value = "景顺长城安泽回报一年持有期混合C"
Set-TextValue $q1.Range("D5") "0.24"
Set-TextValue $q1.Range("E5") "34.82"
Set-TextValue $q1.Range("F5") "0.70"
Set-TextValue $q1.Range("G5") "0.0017"
$q1.Range("H5").Value = 10

# Re-apply the original (untouched) template's cell formatting over the
# whole data range so every row ends up styled consistently (only column A
# carries the bold/bordered "s=2" style) instead of picking up the
# temporary "@" text format used above.
$template.Range("A2:H3").Copy()
$q1.Range("A2:H3").PasteSpecial(-4122)
$template.Range("A2:H3").Copy()
$q1.Range("A4:H5").PasteSpecial(-4122)

# --- 3. Update the "总计" summary sheet -------------------------------------
# Insert a new row above the existing 2021-Q4 summary row, and populate it
# with the 2022-Q1 totals; the old row shifts down and becomes index 1.
$summary.Rows.Item(2).Insert()

$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q1"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 1.03

$summary.Range("A3").Value = 1

# Keep "2021-Q4" as the active/selected sheet, same as before the edit
# (only the sheet list changed - the workbook's active tab did not).
$template.Activate()
